# Update the "取得日時" (retrieved timestamp) column on the "ランサーズ" sheet
# so every data row now reflects the new scrape timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-01-20 02:01:20"

for ($row = 2; $row -le 15; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
